$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values look numeric,
# so Excel keeps them as text (matching the original text-typed column).
$textFormatCells = @(
    "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16",
    "D18", "D19", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D29",
    "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D43",
    "D44", "D45", "D46", "D47", "D49", "D50"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "30.129.74"
$ws.Range("E2").Value = "  +1.00%  "

# Row 3
$ws.Range("D3").Value = "1.892.38"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "0.7385"
$ws.Range("E5").Value = "  -1.35%  "

# Row 6
$ws.Range("D6").Value = "242.84"
$ws.Range("E6").Value = "  +0.25%  "

# Row 7
$ws.Range("E7").Value = "  +0.18%  "

# Row 8
$ws.Range("D8").Value = "0.3172"
$ws.Range("E8").Value = "  +1.37%  "

# Row 9
$ws.Range("D9").Value = "0.07209"
$ws.Range("E9").Value = "  +1.17%  "

# Row 10
$ws.Range("D10").Value = "24.93"
$ws.Range("E10").Value = "  -1.04%  "

# Row 11
$ws.Range("D11").Value = "0.08344"
$ws.Range("E11").Value = "  -1.99%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "2.017.60"
$ws.Range("E12").Value = "  +9.64%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7605"
$ws.Range("E13").Value = "  +0.11%  "

# Row 14
$ws.Range("D14").Value = "5.455"
$ws.Range("E14").Value = "  +1.75%  "

# Row 15
$ws.Range("D15").Value = "93.08"

# Row 16
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "6.166"
$ws.Range("E16").Value = "  +0.09%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.161.51"
$ws.Range("E17").Value = "  +1.11%  "

# Row 18
$ws.Range("D18").Value = "250.42"
$ws.Range("E18").Value = "  +2.89%  "

# Row 19
$ws.Range("D19").Value = "13.64"
$ws.Range("E19").Value = "  -0.34%  "

# Row 20
$ws.Range("D20").Value = "0.000007890"
$ws.Range("E20").Value = "  +1.22%  "

# Row 21
$ws.Range("D21").Value = "2.185.65"
$ws.Range("E21").Value = "  +2.38%  "

# Row 22
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.22%  "

# Row 23
$ws.Range("D23").Value = "7.961"
$ws.Range("E23").Value = "  -0.41%  "

# Row 24
$ws.Range("E24").Value = "  +0.16%  "

# Row 25
$ws.Range("D25").Value = "0.1584"
$ws.Range("E25").Value = "  -0.11%  "

# Row 26
$ws.Range("D26").Value = "9.311"
$ws.Range("E26").Value = "  -0.50%  "

# Row 27
$ws.Range("D27").Value = "164.68"
$ws.Range("E27").Value = "  +1.49%  "

# Row 28
$ws.Range("D28").Value = "18.77"
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("D29").Value = "2.063"
$ws.Range("E29").Value = "  +1.80%  "

# Row 30
$ws.Range("D30").Value = "1.486"
$ws.Range("E30").Value = "  -1.01%  "

# Row 31
$ws.Range("D31").Value = "4.589"
$ws.Range("E31").Value = "  +1.74%  "

# Row 32
$ws.Range("D32").Value = "1.538"
$ws.Range("E32").Value = "  +0.05%  "

# Row 33
$ws.Range("D33").Value = "4.207"
$ws.Range("E33").Value = "  +2.10%  "

# Row 34
$ws.Range("D34").Value = "0.05373"
$ws.Range("E34").Value = "  -0.73%  "

# Row 35
$ws.Range("D35").Value = "1.258"
$ws.Range("E35").Value = "  +1.50%  "

# Row 36
$ws.Range("D36").Value = "0.7781"
$ws.Range("E36").Value = "  +4.22%  "

# Row 37
$ws.Range("E37").Value = "  +0.35%  "

# Row 38
$ws.Range("D38").Value = "2.733"
$ws.Range("E38").Value = "  +0.87%  "

# Row 39
$ws.Range("D39").Value = "0.01970"
$ws.Range("E39").Value = "  +1.42%  "

# Row 40
$ws.Range("E40").Value = "  -0.35%  "

# Row 41
$ws.Range("E41").Value = "  +2.80%  "

# Row 42
$ws.Range("D42").Value = "1.099.09"
$ws.Range("E42").Value = "  +0.80%  "

# Row 43
$ws.Range("D43").Value = "6.085"
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("D44").Value = "72.86"
$ws.Range("E44").Value = "  +0.68%  "

# Row 45
$ws.Range("D45").Value = "0.8715"
$ws.Range("E45").Value = "  +1.85%  "

# Row 46
$ws.Range("D46").Value = "104.51"
$ws.Range("E46").Value = "  +2.14%  "

# Row 47
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.21%  "

# Row 48
$ws.Range("E48").Value = "  +0.40%  "

# Row 49
$ws.Range("D49").Value = "7.603"
$ws.Range("E49").Value = "  -1.17%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.629"
$ws.Range("E50").Value = "  -1.04%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.076.24"
$ws.Range("E51").Value = "  +1.50%  "
